$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.103.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.84%  "

# Row 3
$ws.Range("D3").Value = "'2.378.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.51%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'317.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.30%  "

# Row 6
$ws.Range("D6").Value = "'109.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.32%  "

# Row 7
$ws.Range("E7").Value = "  +0.45%  "

# Row 8
$ws.Range("E8").Value = "  -0.21%  "

# Row 9
$ws.Range("D9").Value = "'0.621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "

# Row 10
$ws.Range("D10").Value = "'41.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.30%  "

# Row 11
$ws.Range("D11").Value = "'0.0933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.36%  "

# Row 12
$ws.Range("D12").Value = "'8.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.37%  "

# Row 13
$ws.Range("E13").Value = "  +0.02%  "

# Row 14
$ws.Range("E14").Value = "  +1.68%  "

# Row 15
$ws.Range("D15").Value = "'16.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "

# Row 16
$ws.Range("D16").Value = "'2.733.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.06%  "

# Row 17
$ws.Range("D17").Value = "'2.374.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "

# Row 18
$ws.Range("D18").Value = "'43.100.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19
$ws.Range("D19").Value = "'7.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "

# Row 20
$ws.Range("E20").Value = "  +1.11%  "

# Row 21
$ws.Range("D21").Value = "'76.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.10%  "

# Row 22
$ws.Range("E22").Value = "  -3.07%  "

# Row 23
$ws.Range("D23").Value = "'270.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "

# Row 24
$ws.Range("E24").Value = "  +1.82%  "

# Row 25
$ws.Range("E25").Value = "  -1.67%  "

# Row 26
$ws.Range("E26").Value = "  +0.37%  "

# Row 27
$ws.Range("E27").Value = "  +0.67%  "

# Row 28
$ws.Range("D28").Value = "'23.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "

# Row 29
$ws.Range("D29").Value = "'2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.24%  "

# Row 30
$ws.Range("D30").Value = "'37.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.63%  "

# Row 31
$ws.Range("D31").Value = "'168.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.36%  "

# Row 32
$ws.Range("D32").Value = "'0.0913"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33
$ws.Range("E33").Value = "  +5.28%  "

# Row 34
$ws.Range("E34").Value = "  -5.44%  "

# Row 35
$ws.Range("E35").Value = "  +16.43%  "

# Row 36
$ws.Range("E36").Value = "  +0.78%  "

# Row 37
$ws.Range("D37").Value = "'4.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.24%  "

# Row 38
$ws.Range("D38").Value = "'0.0363"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "

# Row 39
$ws.Range("E39").Value = "  +0.01%  "

# Row 40
$ws.Range("E40").Value = "  -5.04%  "

# Row 41
$ws.Range("D41").Value = "'105.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.72%  "

# Row 42
$ws.Range("E42").Value = "  +1.25%  "

# Row 43
$ws.Range("E43").Value = "  +5.82%  "

# Row 44
$ws.Range("D44").Value = "'71.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.28%  "

# Row 45
$ws.Range("D45").Value = "'12.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.08%  "

# Row 46
$ws.Range("E46").Value = "  +0.13%  "

# Row 47
$ws.Range("D47").Value = "'115.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.93%  "

# Row 48
$ws.Range("D48").Value = "'81.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +20.23%  "

# Row 49
$ws.Range("E49").Value = "  +3.53%  "

# Row 50
$ws.Range("D50").Value = "'9.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.15%  "

# Row 51
$ws.Range("E51").Value = "  +3.55%  "
